# "ajout force stabilisation dans repere scapula dans OutputFile"
# Rename the two worksheets and leave "Force stabilisation" as the
# active / selected tab with K19 as the active cell.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Feuil1").Name = "Force stabilisation"
$wb.Worksheets.Item("Feuil2").Name = "RoM"

$wsForce = $wb.Worksheets.Item("Force stabilisation")
$wsForce.Activate()
$wsForce.Range("K19").Select()
